# Auto-generated edit script applying F/G column updates across 4 worksheets
# per commit 456a3b4 (gh-pages output regeneration).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1607
$ws.Range("F3").Value = 3356
$ws.Range("F5").Value = 772
$ws.Range("F6").Value = 2400
$ws.Range("F7").Value = 516
$ws.Range("F8").Value = 434
$ws.Range("F9").Value = 261
$ws.Range("F10").Value = 158
$ws.Range("F11").Value = 383
$ws.Range("F12").Value = 1124
$ws.Range("F13").Value = 478
$ws.Range("F14").Value = 235
$ws.Range("F16").Value = 290
$ws.Range("F17").Value = 5016
$ws.Range("F18").Value = 33
$ws.Range("F19").Value = 1410
$ws.Range("F20").Value = 3685
$ws.Range("F22").Value = 237
$ws.Range("F23").Value = 4033
$ws.Range("F24").Value = 5380
$ws.Range("F25").Value = 129
$ws.Range("F26").Value = 993
$ws.Range("F27").Value = 588
$ws.Range("F28").Value = 3443
$ws.Range("F29").Value = 406
$ws.Range("F33").Value = 912
$ws.Range("F34").Value = 1241
$ws.Range("F35").Value = 58
$ws.Range("F36").Value = 75
$ws.Range("F37").Value = 1473
$ws.Range("F38").Value = 156
$ws.Range("F39").Value = 1467
$ws.Range("F40").Value = 58
$ws.Range("F41").Value = 951
$ws.Range("F42").Value = 946
$ws.Range("F43").Value = 540
$ws.Range("F45").Value = 2516
$ws.Range("F46").Value = 96
$ws.Range("F47").Value = 199
$ws.Range("F48").Value = 379
$ws.Range("F49").Value = 3782

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 1054
$ws.Range("F23").Value = 49

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2932

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2932
$ws.Range("F3").Value = 1607
$ws.Range("F4").Value = 3356
$ws.Range("F6").Value = 772
$ws.Range("F7").Value = 2400
$ws.Range("F8").Value = 516
$ws.Range("F9").Value = 434
$ws.Range("F10").Value = 261
$ws.Range("F11").Value = 1054
$ws.Range("F12").Value = 158
$ws.Range("F13").Value = 383
$ws.Range("F14").Value = 1124
$ws.Range("F15").Value = 478
$ws.Range("F16").Value = 235
$ws.Range("F18").Value = 290
$ws.Range("F19").Value = 5016
$ws.Range("F20").Value = 1410
$ws.Range("F21").Value = 4035
$ws.Range("F22").Value = 5381
$ws.Range("F23").Value = 129
$ws.Range("F24").Value = 993
$ws.Range("F25").Value = 588
$ws.Range("F26").Value = 3443
$ws.Range("F27").Value = 406
$ws.Range("F31").Value = 912
$ws.Range("F32").Value = 1241
$ws.Range("F33").Value = 58
$ws.Range("F34").Value = 75
$ws.Range("F35").Value = 1473
$ws.Range("F36").Value = 156
$ws.Range("F37").Value = 1467
$ws.Range("F39").Value = 951
$ws.Range("F41").Value = 540
$ws.Range("F44").Value = 49
$ws.Range("F45").Value = 2516
$ws.Range("F46").Value = 96
$ws.Range("F47").Value = 199
$ws.Range("F48").Value = 379
$ws.Range("F49").Value = 3782

# Special case: G21 on 展览 becomes sold-out text, and F21 updates too
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F21").Value = 236
$ws.Range("G21").Value = "已售罄"
